# EPBDS-12620 Difference in error response structure between kafka and webservice call
#
# The rule table's last row ("R40") returned the literal greeting "Good Night".
# It is replaced with the text '= error("fail")' so the test rule intentionally
# fails/returns an error value. The leading "=" must stay literal text (not be
# interpreted as a formula), which is why it is entered with a leading
# apostrophe - Excel then stores it as a shared string and marks the cell's
# style with quotePrefix so the apostrophe isn't shown/re-parsed as a formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E11").Value = "'= error(""fail"")"

# Mirror the active selection recorded in the saved view state.
[void]$ws.Range("F9").Select()
